# Misc updates related to Reach Ranking, etc.
#
# The "PROTECTION" reach-ranking sheet gets an extra Entiat River reach
# inserted, several scores recomputed for the surrounding rows, and the
# Twisp River section gets renumbered/consolidated (three old "Twisp
# River Lower/Middle" rows collapse into two "Twisp River Middle" rows
# plus one updated "Twisp River Middle 06" row), which shrinks the used
# range from A1:W13 down to A1:W12.
#
# Rows 2 (Big Meadow Creek 01) and 3 (Entiat River Potato 05) are left
# untouched. Rows 4-12 are (re)written with their final values, and the
# former row 13 is deleted outright (which also fixes up the sheet's
# dimension to A1:W12 automatically).

$data = @(
    @('Entiat River Potato 06', 'Entiat', 'Entiat River-Potato Creek', 'yes', 'yes', 'yes', 5, 5, 5, 5, 5, 3, 3, 3, 5, 5, 3, 4, 3, 36, 0.8, 3, 3),
    @('Methow River Fawn 04', 'Methow', 'Methow River-Fawn Creek', 'yes', 'yes', 'yes', 5, 5, 5, 5, 3, 3, 5, 1, 5, 5, 3, 4, 1, 32, 0.7111111111111111, 5, 3),
    @('Methow River Thompson 08', 'Methow', 'Methow River-Thompson Creek', 'yes', 'yes', 'yes', 3, 5, 4, 5, 5, 1, 3, 3, 5, 1, 1, 1, 5, 32, 0.7111111111111111, 5, 3),
    @('Nason Creek Lower 01', 'Wenatchee', 'Lower Nason Creek', 'yes', 'yes', 'yes', 3, 5, 4, 5, 5, 5, 3, 5, 5, 3, 3, 3, 1, 36, 0.8, 3, 3),
    @('Nason Creek Lower 02', 'Wenatchee', 'Lower Nason Creek', 'yes', 'yes', 'yes', 3, 5, 4, 5, 5, 5, 3, 1, 5, 5, 5, 5, 3, 36, 0.8, 3, 3),
    @('Nason Creek Lower 03', 'Wenatchee', 'Lower Nason Creek', 'yes', 'yes', 'yes', 3, 5, 4, 5, 5, 3, 3, 5, 5, 3, 3, 3, 1, 34, 0.7555555555555555, 5, 3),
    @('Twisp River Middle 01', 'Methow', 'Middle Twisp River', 'yes', 'yes', 'yes', 3, 3, 3, 5, 5, 3, 3, 5, 5, 3, 3, 3, 1, 33, 0.7333333333333333, 5, 3),
    @('Twisp River Middle 02', 'Methow', 'Middle Twisp River', 'yes', 'yes', 'yes', 3, 3, 3, 5, 5, 5, 3, 5, 5, 3, 3, 3, 1, 35, 0.7777777777777778, 5, 3),
    @('Twisp River Middle 06', 'Methow', 'Middle Twisp River', 'yes', 'yes', 'yes', 5, 5, 5, 5, 1, 5, 5, 3, 5, 5, 1, 3, 3, 35, 0.7777777777777778, 5, 3)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowVals = $data[$i]
    for ($col = 1; $col -le $rowVals.Count; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowVals[$col - 1]
    }
}

# The old trailing row (previously row 13, "Twisp River Middle 05") is no
# longer present in the updated table; delete it so the used range shrinks
# from A1:W13 to A1:W12.
$ws.Rows.Item(13).Delete()
